# Generate Report for Handoff
#
# This updates the localization-status report for the four "Ready for
# handoff" rows (1debd559-..., 24ae5118-..., 30a66ade-..., a0378fe7-...)
# which live in rows 4-7 of the "zh-cn" and "de-de" worksheets:
#   - Priority moves from "low" to "ht" (handoff triggered)
#   - The handoff timestamps are refreshed to the new generation time
#
# zh-cn: Latest Handoff Datetime (column H) -> 2016-08-17 06:27:26
# de-de: Latest HO Xliff Generate Date (shared with Overview column G,
#         and de-de's own Latest Handoff Datetime) -> 2016-08-17 06:27:31

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7 -> Priority (E) and Latest Handoff Datetime (H)
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-17 06:27:26"
}

# de-de: rows 4-7 -> Priority (E) and Latest Handoff Datetime (H)
# (the de-de handoff datetime / HO xliff generate date text is the same
# value shown by the Overview sheet, so it is refreshed there too)
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-17 06:27:31"
}

# Overview: rows 4-7 -> Latest HO Xliff Generate Date (G)
for ($r = 4; $r -le 7; $r++) {
    $overview.Cells.Item($r, 7).Value = "2016-08-17 06:27:31"
}
